$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and reordered rows 32-35)
# Leading apostrophes force Excel to treat the assigned values as text,
# matching the original inlineStr/shared-string cell type and preventing
# numeric-looking strings (e.g. "6.750", "1.040") from being coerced to numbers.
$ws.Range("D2").Value = "'27.825.18"
$ws.Range("E2").Value = "'  +0.66%  "
$ws.Range("D3").Value = "'1.859.73"
$ws.Range("E3").Value = "'  +0.37%  "
$ws.Range("E4").Value = "'  +0.39%  "
$ws.Range("D5").Value = "'323.75"
$ws.Range("E5").Value = "'  +0.70%  "
$ws.Range("E6").Value = "'  +0.26%  "
$ws.Range("D7").Value = "'0.4418"
$ws.Range("D8").Value = "'0.3825"
$ws.Range("E8").Value = "'  +1.80%  "
$ws.Range("D9").Value = "'0.07446"
$ws.Range("E9").Value = "'  +0.43%  "
$ws.Range("D10").Value = "'0.8889"
$ws.Range("E10").Value = "'  +1.40%  "
$ws.Range("D11").Value = "'21.66"
$ws.Range("E11").Value = "'  +0.79%  "
$ws.Range("D12").Value = "'1.874.08"
$ws.Range("E12").Value = "'  +0.95%  "
$ws.Range("D13").Value = "'5.542"
$ws.Range("E13").Value = "'  +0.41%  "
$ws.Range("D14").Value = "'6.750"
$ws.Range("E14").Value = "'  +0.64%  "
$ws.Range("D15").Value = "'0.07214"
$ws.Range("D16").Value = "'86.27"
$ws.Range("E16").Value = "'  +3.98%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("E17").Value = "'  +0.38%  "
$ws.Range("D18").Value = "'0.000009113"
$ws.Range("E18").Value = "'  +0.69%  "
$ws.Range("E19").Value = "'  +0.33%  "
$ws.Range("D20").Value = "'15.60"
$ws.Range("E20").Value = "'  +0.93%  "
$ws.Range("D21").Value = "'27.841.45"
$ws.Range("E21").Value = "'  +0.70%  "
$ws.Range("D22").Value = "'5.305"
$ws.Range("E22").Value = "'  +0.64%  "
$ws.Range("E23").Value = "'  +0.36%  "
$ws.Range("D24").Value = "'2.102.50"
$ws.Range("E24").Value = "'  +1.46%  "
$ws.Range("D25").Value = "'2.063"
$ws.Range("E25").Value = "'  +6.07%  "
$ws.Range("D26").Value = "'159.35"
$ws.Range("E26").Value = "'  +1.07%  "
$ws.Range("D27").Value = "'18.80"
$ws.Range("E27").Value = "'  +0.14%  "
$ws.Range("D28").Value = "'2.014"
$ws.Range("E28").Value = "'  +3.98%  "
$ws.Range("D29").Value = "'5.387"
$ws.Range("E29").Value = "'  +1.53%  "
$ws.Range("D30").Value = "'118.82"
$ws.Range("E30").Value = "'  +2.00%  "
$ws.Range("D31").Value = "'0.09116"
$ws.Range("E31").Value = "'  +0.38%  "
$ws.Range("B32").Value = "'ARBITRUM"
$ws.Range("C32").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.219"
$ws.Range("E32").Value = "'  +0.92%  "
$ws.Range("B33").Value = "'ImmutableX"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7766"
$ws.Range("E33").Value = "'  +1.11%  "
$ws.Range("B34").Value = "'Filecoin"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.621"
$ws.Range("E34").Value = "'  +2.38%  "
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.024"
$ws.Range("E35").Value = "'  +4.77%  "
$ws.Range("E36").Value = "'  +0.38%  "
$ws.Range("D37").Value = "'1.159"
$ws.Range("E37").Value = "'  +0.19%  "
$ws.Range("D38").Value = "'0.01988"
$ws.Range("E38").Value = "'  +0.37%  "
$ws.Range("D39").Value = "'0.05332"
$ws.Range("E39").Value = "'  +0.93%  "
$ws.Range("D40").Value = "'2.874"
$ws.Range("E40").Value = "'  +2.06%  "
$ws.Range("D41").Value = "'0.5228"
$ws.Range("E41").Value = "'  +0.88%  "
$ws.Range("D42").Value = "'6.971"
$ws.Range("E42").Value = "'  +3.72%  "
$ws.Range("D43").Value = "'0.1682"
$ws.Range("E43").Value = "'  +0.40%  "
$ws.Range("D44").Value = "'8.827"
$ws.Range("E44").Value = "'  +2.86%  "
$ws.Range("D45").Value = "'111.21"
$ws.Range("E45").Value = "'  +2.06%  "
$ws.Range("D46").Value = "'10.80"
$ws.Range("E46").Value = "'  +1.98%  "
$ws.Range("D47").Value = "'1.037"
$ws.Range("E47").Value = "'  +0.37%  "
$ws.Range("D48").Value = "'0.06593"
$ws.Range("E48").Value = "'  +3.13%  "
$ws.Range("D49").Value = "'1.721"
$ws.Range("E49").Value = "'  -0.12%  "
$ws.Range("D50").Value = "'0.4740"
$ws.Range("E50").Value = "'  +1.68%  "
$ws.Range("D51").Value = "'1.895"
$ws.Range("E51").Value = "'  +0.46%  "